# Insert a new data row at row 51 (pushing the existing rows 51-85 down to
# 52-86) and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 51..85 down to 52..86, new blank row appears at 51
# (format is inherited from the row above, matching row D's date style).
$ws.Rows.Item(51).Insert()

$ws.Range("A51").Value = 11
$ws.Range("B51").Value = "Vega Monumental Concepción"
$ws.Range("C51").Value = "Bíobío"
$ws.Range("D51").Value = 44680
$ws.Range("E51").Value = 8
$ws.Range("F51").Value = 100112001
$ws.Range("G51").Value = "Berenjena"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 6000
$ws.Range("L51").Value = 7000
$ws.Range("M51").Value = 6500
$ws.Range("N51").Value = "$/caja 60 unidades"
$ws.Range("O51").Value = "Región Metropolitana"
$ws.Range("P51").Value = 108
$ws.Range("Q51").Value = 60
$ws.Range("R51").Value = "Hortaliza"
